$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.327.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.665.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.93%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.78%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5353"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.67%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.88%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2662"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06408"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07852"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.562"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.643.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.892.95"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5538"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8187"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.69%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.349.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.10%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.684"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.37%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.29"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.038"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.65"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.80%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1232"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.215"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05870"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.83%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.650"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.283"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.610"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9689"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.72%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.420"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5833"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.77%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01606"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.01%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8690"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.07%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.848"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.80%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.32"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.61%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.052.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.30%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.89%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.803.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.57%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.07%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.014"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.17%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.020"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.27%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.49%  "
